# Applies the "nuevo periodo 2509" update to the Estado de Cuenta worksheet:
#  - Adds a new trailing data row for period 2509 (same worker/value pattern
#    as the existing rows) right before the old last row.
#  - The old last row (2508) becomes a "normal" interior row, and the new
#    row (2509) takes over the special bottom-border ("last row") styling.
#  - Bumps the summary counters (Cant. Periodos, Valor Mora) accordingly.
#  - Centers the "Periodo Mora" column for every data row.
#  - Leaves the signature block (underline + labels) intact; it is simply
#    pushed down a row by the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats
$xlCenter = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

# --- 1. Insert a new row right after the current last data row (41), ---
#        pushing the blank spacer + signature block rows down by one.
$ws.Rows.Item(42).Insert()

# --- 2. The inserted row 42 should look like the old "last row" (41), ---
#        which still carries the special bottom-border formatting.
$ws.Range("B41:J41").Copy()
$ws.Range("B42:J42").PasteSpecial($xlPasteFormats)

# --- 3. Row 41 goes back to being a normal interior row (format like 40). ---
$ws.Range("B40:J40").Copy()
$ws.Range("B41:J41").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- 4. Fill in the new period-2509 row with the same worker data. ---
$ws.Cells.Item(42, 2).Value2 = "CC"
$ws.Cells.Item(42, 3).Value2 = "1104425093"
$ws.Cells.Item(42, 4).Value2 = "LULIETH PATRICIA MORALES PRASCA"
$ws.Cells.Item(42, 5).Value2 = "2509"
$ws.Cells.Item(42, 6).Value2 = 46400
$ws.Cells.Item(42, 7).Value2 = 1160000

# --- 5. Center the "Periodo Mora" column for every data row (16-42). ---
$ws.Range("E16:E42").HorizontalAlignment = $xlCenter

# --- 6. Update the summary counters. ---
$ws.Range("F13").Value2 = 27
$ws.Range("E11").Value2 = 1252800

Write-Output "Edit applied"
